{"js": "// Replace the date line and the twenty-five \"NNN\u00f7N=\" problems in the table\n// with the values from the new day's worksheet. Every source string in this\n// document is unique, so a simple search-and-replace per pair is safe and\n// unambiguous.\nconst replacements = [\n  [\"2024-08-14 Wednesday\", \"2024-08-15 Thursday\"],\n  [\"959\u00f79=\", \"132\u00f75=\"],\n  [\"955\u00f72=\", \"247\u00f78=\"],\n  [\"438\u00f72=\", \"693\u00f79=\"],\n  [\"333\u00f74=\", \"868\u00f76=\"],\n  [\"612\u00f79=\", \"244\u00f76=\"],\n  [\"907\u00f74=\", \"666\u00f74=\"],\n  [\"950\u00f72=\", \"507\u00f76=\"],\n  [\"440\u00f73=\", \"772\u00f79=\"],\n  [\"722\u00f78=\", \"515\u00f78=\"],\n  [\"603\u00f77=\", \"705\u00f75=\"],\n  [\"927\u00f75=\", \"123\u00f78=\"],\n  [\"583\u00f72=\", \"491\u00f73=\"],\n  [\"703\u00f74=\", \"690\u00f74=\"],\n  [\"236\u00f78=\", \"468\u00f75=\"],\n  [\"988\u00f78=\", \"629\u00f78=\"],\n  [\"265\u00f75=\", \"423\u00f75=\"],\n  [\"760\u00f79=\", \"695\u00f72=\"],\n  [\"421\u00f77=\", \"193\u00f73=\"],\n  [\"755\u00f77=\", \"420\u00f77=\"],\n  [\"218\u00f79=\", \"328\u00f73=\"],\n  [\"387\u00f79=\", \"387\u00f73=\"],\n  [\"483\u00f79=\", \"153\u00f78=\"],\n  [\"927\u00f72=\", \"472\u00f78=\"],\n  [\"280\u00f78=\", \"549\u00f74=\"],\n  [\"317\u00f76=\", \"439\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the twenty-five \"NNN\u00f7N=\" problems in the table\n# with the values from the new day's worksheet. Every source string in this\n# document is unique, so a simple Find/Replace per pair is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-08-14 Wednesday\", \"2024-08-15 Thursday\"),\n  @(\"959\u00f79=\", \"132\u00f75=\"),\n  @(\"955\u00f72=\", \"247\u00f78=\"),\n  @(\"438\u00f72=\", \"693\u00f79=\"),\n  @(\"333\u00f74=\", \"868\u00f76=\"),\n  @(\"612\u00f79=\", \"244\u00f76=\"),\n  @(\"907\u00f74=\", \"666\u00f74=\"),\n  @(\"950\u00f72=\", \"507\u00f76=\"),\n  @(\"440\u00f73=\", \"772\u00f79=\"),\n  @(\"722\u00f78=\", \"515\u00f78=\"),\n  @(\"603\u00f77=\", \"705\u00f75=\"),\n  @(\"927\u00f75=\", \"123\u00f78=\"),\n  @(\"583\u00f72=\", \"491\u00f73=\"),\n  @(\"703\u00f74=\", \"690\u00f74=\"),\n  @(\"236\u00f78=\", \"468\u00f75=\"),\n  @(\"988\u00f78=\", \"629\u00f78=\"),\n  @(\"265\u00f75=\", \"423\u00f75=\"),\n  @(\"760\u00f79=\", \"695\u00f72=\"),\n  @(\"421\u00f77=\", \"193\u00f73=\"),\n  @(\"755\u00f77=\", \"420\u00f77=\"),\n  @(\"218\u00f79=\", \"328\u00f73=\"),\n  @(\"387\u00f79=\", \"387\u00f73=\"),\n  @(\"483\u00f79=\", \"153\u00f78=\"),\n  @(\"927\u00f72=\", \"472\u00f78=\"),\n  @(\"280\u00f78=\", \"549\u00f74=\"),\n  @(\"317\u00f76=\", \"439\u00f77=\")\n)\n\nforeach ($pair in $pairs) {\n  $from = $pair[0]\n  $to = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $from\n  $find.Replacement.Text = $to\n  $find.Execute($from, $false, $false, $false, $false, $false, $true, 1, $false, $to, 2) | Out-Null\n}\n"}
